$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$txt = "Microsoft Office User:Marc" + [char]10 + "include scenes" + [char]10 + "areaX:scene,scene,scene;areaY:scene,scene"
$cm = $ws.Range("F90").AddComment($txt)
$cm2 = $ws.Range("F91").AddComment($txt)
